# Applies the edit described by the diff:
#  - Removes three rows (original row numbers 13, 27, 31) that correspond to
#    entries whose DIVISION is a city not part of the standard province list
#    for this filtered region file (Kabankalan City, Cadiz City, Himamaylan City).
#  - Adjusts the width of columns C, F and Y.
# Deleting rows automatically shifts everything below up and keeps the
# dataValidation sqref / sheet dimension in sync, matching the target XML.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Delete rows from bottom to top so the remaining row indices stay valid.
$ws.Rows.Item(31).Delete()
$ws.Rows.Item(27).Delete()
$ws.Rows.Item(13).Delete()

# Column width adjustments (XML "width" in characters = ColumnWidth + 0.83
# for this workbook's default font, based on empirical measurement).
$ws.Columns.Item(3).ColumnWidth = 9.17    # column C: 17 -> 10
$ws.Columns.Item(6).ColumnWidth = 15.17   # column F: 20 -> 16
$ws.Columns.Item(25).ColumnWidth = 119.17 # column Y: 170 -> 120
